# "mend the equip effect"
# The ItemGift "I组合礼包2" style gift list (cell B4, shared string used by
# row 4 / item 22031001) erroneously included two equip-type entries
# (2;21200101;... and 2;21400101;...) appended after the intended item
# list. Remove them so the reward list only grants the intended items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemGift")

$ws.Range("B4").Value = "1;22033001;1;100;5|1;22033002;1;100;5|1;22032007;1;100;1|1;22033013;1;100;5|1;22033014;1;100;3|1;22033015;1;100;3"

# Re-saving the workbook (as happened in the authoring environment) also
# touched a couple of cosmetic/presentation settings; reproduce the ones
# that are reachable through the Excel object model.

# Sheet background ("window"/lt1) theme colour.
$themeColors = $wb.Theme.ThemeColorScheme
$lt1 = $themeColors.Colors(2)
$lt1.RGB = 13494986   # RGB(202,234,205) == CAEACD

# Page setup for the worksheet (paper size / orientation).
$pageSetup = $ws.PageSetup
$pageSetup.Orientation = 1
$pageSetup.PaperSize = 9
